$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.440.58"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.875.93"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.05"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.79"
$ws.Range("E8").Value = "  +7.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.332"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("D12").Value = "2.149.82"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "1.866.89"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.685"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.74"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "35.421.83"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.14"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.65"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.34"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.27"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").Value = "  +4.96%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.89"
$ws.Range("E27").Value = "  +23.49%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("B31").Value = "BinanceUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.01"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.05"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.05"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +21.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.832"
$ws.Range("E35").Value = "  +18.48%  "
$ws.Range("E36").Value = "  +5.92%  "
$ws.Range("E37").Value = "  +7.24%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0204"
$ws.Range("E39").Value = "  +4.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.92"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "1.352.12"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.31"
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0599"
$ws.Range("E43").Value = "  +14.44%  "
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.02"
$ws.Range("E44").Value = "  +55.91%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.35"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.65"
$ws.Range("E47").Value = "  +6.55%  "
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").Value = "2.058.37"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.45"
$ws.Range("E51").Value = "  +1.53%  "
